# Updated cryptos list values (price + 1h volume change) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells that would otherwise be auto-parsed as numbers (single
# decimal point) need an explicit text format applied first so Excel keeps them
# as literal strings, matching the source data which stores prices as text.
$textPriceRows = @(5, 6, 9, 10, 11, 12, 17, 23, 24, 27, 28, 30, 31, 32, 33, 38, 41, 43, 44, 49, 50, 51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "48.084.60"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.502.65"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "320.67"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "107.36"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "39.52"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "20.15"
$ws.Range("E11").Value = "  +7.37%  "
$ws.Range("D12").Value = "0.0811"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "2.894.22"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "2.503.50"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.835"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "47.945.54"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "0.0₃0937"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").Value = "278.59"
$ws.Range("E23").Value = "  +12.49%  "
$ws.Range("D24").Value = "71.45"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "25.85"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "9.65"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").Value = "35.24"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  -4.72%  "
$ws.Range("D32").Value = "49.51"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "121.06"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "21.24"
$ws.Range("E43").Value = "  -5.83%  "
$ws.Range("D44").Value = "0.0300"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "2.013.44"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "8.99"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "5.18"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "80.07"
$ws.Range("E51").Value = "  +2.82%  "
